$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows 3-6 hold the "想去人数" (want-to-go count) values in column F
$wsExhibitions = $wb.Worksheets.Item("展览")
$wsExhibitions.Range("F3").Value = 2101
$wsExhibitions.Range("F4").Value = 868
$wsExhibitions.Range("F5").Value = 1295
$wsExhibitions.Range("F6").Value = 363

# Sheet "全部类型" (All Types) - same underlying data, but at different rows (3, 6, 7, 8)
$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F3").Value = 2101
$wsAllTypes.Range("F6").Value = 868
$wsAllTypes.Range("F7").Value = 1295
$wsAllTypes.Range("F8").Value = 363
